$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Team Brian Drumm (row 3) correction: General Knowledge score 5 -> 4
$ws.Range("E3").Value = 4

# Team Brian Noonan (row 8) added/corrected scores
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 0.5
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 3.5
$ws.Range("G8").Value = 4

$excel.CalculateFullRebuild()

# Update the active selection to I3
$ws.Activate()
$ws.Range("I3").Select()

$wb.Save()
